$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (46074 -> 46075) for every data row (rows 2 through 435).
$firstRow = 2
$lastRow = 435
$col = 3  # Column C

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Value2 -eq 46074) {
        $cell.Value2 = 46075
    }
}
